$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 325
$ws.Range("I11").Value = 325
$ws.Range("K11").Value = 325
$ws.Range("M11").Value = -185
$ws.Range("H32").Value = 83338830
$ws.Range("J32").Value = 7500
$ws.Range("L32").Value = 7500
$ws.Range("N32").Value = -8152
$ws.Range("H88").Value = 13849.4
$ws.Range("I88").Value = 11500
$ws.Range("K88").Value = 11500
$ws.Range("M88").Value = -11094
$ws.Range("H91").Value = 13849.4
$ws.Range("I91").Value = 11500
$ws.Range("K91").Value = 11500
$ws.Range("M91").Value = -10096
$ws.Range("H98").Value = 658.1739
$ws.Range("I98").Value = 651.7273
$ws.Range("K98").Value = 651.7273
$ws.Range("M98").Value = 846.2727
$ws.Range("H114").Value = 29250
$ws.Range("J114").Value = 29250
$ws.Range("L114").Value = 29250
$ws.Range("N114").Value = -37928
$ws.Range("H122").Value = 658.1739
$ws.Range("I122").Value = 651.7273
$ws.Range("K122").Value = 1955.1819
$ws.Range("M122").Value = 494.8181
$ws.Range("H137").Value = 16668787
$ws.Range("I137").Value = 1860
$ws.Range("K137").Value = 5580
$ws.Range("M137").Value = -3030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17127.924
$ws.Range("I32").Value = 16595.783
$ws.Range("K32").Value = 16595.783
$ws.Range("M32").Value = -16308.783
$ws.Range("H57").Value = 8000
$ws.Range("I57").Value = 8000
$ws.Range("K57").Value = 8000
$ws.Range("M57").Value = -7516
$ws.Range("H74").Value = 1773.1
$ws.Range("I74").Value = 1404
$ws.Range("K74").Value = 1404
$ws.Range("M74").Value = -530
$ws.Range("H77").Value = 1773.1
$ws.Range("I77").Value = 1404
$ws.Range("K77").Value = 7020
$ws.Range("M77").Value = -2652
$ws.Range("H88").Value = 44259.25
$ws.Range("I88").Value = 678
$ws.Range("K88").Value = 678
$ws.Range("M88").Value = -272
$ws.Range("H91").Value = 44259.25
$ws.Range("I91").Value = 678
$ws.Range("K91").Value = 678
$ws.Range("M91").Value = 726
$ws.Range("H97").Value = 1393.9048
$ws.Range("I97").Value = 1372.6428
$ws.Range("K97").Value = 1372.6428
$ws.Range("M97").Value = -876.6428000000001
$ws.Range("H126").Value = 2040
$ws.Range("I126").Value = 2040
$ws.Range("K126").Value = 6120
$ws.Range("M126").Value = -3650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 599.25
$ws.Range("I11").Value = 549
$ws.Range("J11").Value = 750
$ws.Range("K11").Value = 549
$ws.Range("L11").Value = 750
$ws.Range("M11").Value = -409
$ws.Range("N11").Value = -1030
$ws.Range("H86").Value = 2573.9583
$ws.Range("I86").Value = 1915.2667
$ws.Range("J86").Value = 3671.7778
$ws.Range("K86").Value = 1915.2667
$ws.Range("L86").Value = 3671.7778
$ws.Range("M86").Value = -792.2666999999999
$ws.Range("N86").Value = -5917.7778
$ws.Range("H89").Value = 2573.9583
$ws.Range("I89").Value = 1915.2667
$ws.Range("J89").Value = 3671.7778
$ws.Range("K89").Value = 9576.333499999999
$ws.Range("L89").Value = 18358.889
$ws.Range("M89").Value = -3960.333499999999
$ws.Range("N89").Value = -29590.889
$ws.Range("H99").Value = 2148.6667
$ws.Range("I99").Value = 1693.5714
$ws.Range("K99").Value = 1693.5714
$ws.Range("M99").Value = -195.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35716772
$ws.Range("I31").Value = 47620604
$ws.Range("J31").Value = 5282
$ws.Range("K31").Value = 47620604
$ws.Range("L31").Value = 5282
$ws.Range("M31").Value = -47620309
$ws.Range("N31").Value = -5872
$ws.Range("H34").Value = 35716772
$ws.Range("I34").Value = 47620604
$ws.Range("J34").Value = 5282
$ws.Range("K34").Value = 47620604
$ws.Range("L34").Value = 5282
$ws.Range("M34").Value = -47620402
$ws.Range("N34").Value = -5686
$ws.Range("H62").Value = 26599.357
$ws.Range("I62").Value = 4979
$ws.Range("K62").Value = 4979
$ws.Range("M62").Value = -4355
$ws.Range("H65").Value = 26599.357
$ws.Range("I65").Value = 4979
$ws.Range("K65").Value = 24895
$ws.Range("M65").Value = -21775
$ws.Range("H132").Value = 33344704
$ws.Range("I132").Value = 47623388
$ws.Range("J132").Value = 27777.25
$ws.Range("K132").Value = 142870164
$ws.Range("L132").Value = 83331.75
$ws.Range("M132").Value = -142867634
$ws.Range("N132").Value = -88391.75
$ws.Range("H134").Value = 2713.6155
$ws.Range("J134").Value = 3315.875
$ws.Range("L134").Value = 9947.625
$ws.Range("N134").Value = -15017.625
$ws.Range("H141").Value = 107597.05
$ws.Range("J141").Value = 110020.71
$ws.Range("L141").Value = 110020.71
$ws.Range("N141").Value = -120380.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 149
$ws.Range("I2").Value = 14.25
$ws.Range("J2").Value = 283.75
$ws.Range("K2").Value = 85.5
$ws.Range("L2").Value = 1702.5
$ws.Range("M2").Value = 27.5
$ws.Range("N2").Value = -1928.5
$ws.Range("H5").Value = 1172
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 600
$ws.Range("M5").Value = -488
$ws.Range("H23").Value = 1801.3846
$ws.Range("J23").Value = 2124.625
$ws.Range("L23").Value = 6373.875
$ws.Range("N23").Value = -6843.875
$ws.Range("H37").Value = 142951420
$ws.Range("J37").Value = 142951420
$ws.Range("L37").Value = 428854260
$ws.Range("N37").Value = -428854484
$ws.Range("H38").Value = 160.63637
$ws.Range("I38").Value = 123.71429
$ws.Range("J38").Value = 225.25
$ws.Range("K38").Value = 371.14287
$ws.Range("L38").Value = 675.75
$ws.Range("M38").Value = -24.14287000000002
$ws.Range("N38").Value = -1369.75
$ws.Range("H135").Value = 1172
$ws.Range("I135").Value = 200
$ws.Range("K135").Value = 1800
$ws.Range("M135").Value = 735

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 46665.832
$ws.Range("J52").Value = 46665.832
$ws.Range("L52").Value = 46665.832
$ws.Range("N52").Value = -47183.832
$ws.Range("H70").Value = 7574
$ws.Range("I70").Value = 7446.8335
$ws.Range("J70").Value = 7792
$ws.Range("K70").Value = 7446.8335
$ws.Range("L70").Value = 7792
$ws.Range("M70").Value = -7176.8335
$ws.Range("N70").Value = -8332
$ws.Range("H73").Value = 7574
$ws.Range("I73").Value = 7446.8335
$ws.Range("J73").Value = 7792
$ws.Range("K73").Value = 7446.8335
$ws.Range("L73").Value = 7792
$ws.Range("M73").Value = -6510.8335
$ws.Range("N73").Value = -9664
$ws.Range("H107").Value = 954.8889
$ws.Range("I107").Value = 900
$ws.Range("K107").Value = 900
$ws.Range("M107").Value = 1020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3777.6028
$ws.Range("J7").Value = 5083.4
$ws.Range("L7").Value = 5083.4
$ws.Range("N7").Value = -5307.4
$ws.Range("H16").Value = 4780.2
$ws.Range("I16").Value = 3633
$ws.Range("J16").Value = 6501
$ws.Range("K16").Value = 3633
$ws.Range("L16").Value = 6501
$ws.Range("M16").Value = -3463
$ws.Range("N16").Value = -6841
$ws.Range("H68").Value = 2972.4473
$ws.Range("I68").Value = 2527.1738
$ws.Range("J68").Value = 3655.2
$ws.Range("K68").Value = 2527.1738
$ws.Range("L68").Value = 3655.2
$ws.Range("M68").Value = -1778.1738
$ws.Range("N68").Value = -5153.2
$ws.Range("H71").Value = 2972.4473
$ws.Range("I71").Value = 2527.1738
$ws.Range("J71").Value = 3655.2
$ws.Range("K71").Value = 12635.869
$ws.Range("L71").Value = 18276
$ws.Range("M71").Value = -8891.869000000001
$ws.Range("N71").Value = -25764
$ws.Range("H74").Value = 88000
$ws.Range("J74").Value = 88000
$ws.Range("L74").Value = 88000
$ws.Range("N74").Value = -89996
$ws.Range("H77").Value = 88000
$ws.Range("J77").Value = 88000
$ws.Range("L77").Value = 264000
$ws.Range("N77").Value = -273984
$ws.Range("H100").Value = 1998
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 3777.6028
$ws.Range("J126").Value = 5083.4
$ws.Range("L126").Value = 15250.2
$ws.Range("N126").Value = -20190.2
$ws.Range("H136").Value = 3623.558
$ws.Range("I136").Value = 2479.5862
$ws.Range("K136").Value = 7438.758600000001
$ws.Range("M136").Value = -4888.758600000001
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 27614.8
$ws.Range("I32").Value = 27614.8
$ws.Range("K32").Value = 27614.8
$ws.Range("M32").Value = -27297.8
$ws.Range("H130").Value = 64443.5
$ws.Range("J130").Value = 64443.5
$ws.Range("L130").Value = 64443.5
$ws.Range("N130").Value = -74483.5
$ws.Range("H132").Value = 2997.8333
$ws.Range("I132").Value = 2486.5
$ws.Range("K132").Value = 7459.5
$ws.Range("M132").Value = -4929.5
$ws.Range("H140").Value = 45503.285
$ws.Range("J140").Value = 45503.285
$ws.Range("L140").Value = 45503.285
$ws.Range("N140").Value = -55863.285
